$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new accelerometer sample is prepended to the dataset (becoming the new row 2),
# pushing every existing data row down by one. We shift values directly instead of
# using Rows.Insert() so the plain (unstyled) formatting of the data rows is preserved
# — Insert() would otherwise copy the header row's bold/centered style onto the new row.
for ($r = 21; $r -ge 2; $r--) {
    $x = $ws.Cells.Item($r, 1).Value()
    $y = $ws.Cells.Item($r, 2).Value()
    $z = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($r + 1, 1).Value = $x
    $ws.Cells.Item($r + 1, 2).Value = $y
    $ws.Cells.Item($r + 1, 3).Value = $z
}

# Write the new first sample into row 2.
$ws.Cells.Item(2, 1).Value = -2.568546533584593
$ws.Cells.Item(2, 2).Value = 8.544089555740355
$ws.Cells.Item(2, 3).Value = -1.202380612492562

# The shift above duplicated the former last row (old row 21) into row 22; the dataset
# only spans down to row 21 now, so clear that stale trailing row.
$ws.Range("A22:C22").Clear()
